$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''56.122.57'
$ws.Range('E2').Value = '  +3.40%  '
$ws.Range('D3').Value = '''2.314.11'
$ws.Range('E3').Value = '  +2.09%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''517.82'
$ws.Range('E5').Value = '  +4.36%  '
$ws.Range('D6').Value = '''133.25'
$ws.Range('E6').Value = '  +3.37%  '
$ws.Range('D7').Value = '''0.994'
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('D8').Value = '''0.535'
$ws.Range('E8').Value = '  +1.60%  '
$ws.Range('D9').Value = '''2.337.02'
$ws.Range('E9').Value = '  +2.82%  '
$ws.Range('D10').Value = '''0.104'
$ws.Range('E10').Value = '  +8.69%  '
$ws.Range('D11').Value = '''0.155'
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('D12').Value = '''5.16'
$ws.Range('E12').Value = '  +7.68%  '
$ws.Range('D13').Value = '''0.343'
$ws.Range('E13').Value = '  +1.24%  '
$ws.Range('D14').Value = '''24.03'
$ws.Range('E14').Value = '  +4.91%  '
$ws.Range('D15').Value = '''2.730.01'
$ws.Range('E15').Value = '  +2.29%  '
$ws.Range('D16').Value = '''56.278.41'
$ws.Range('E16').Value = '  +3.74%  '
$ws.Range('D17').Value = '''0.0000136'
$ws.Range('E17').Value = '  +4.57%  '
$ws.Range('D18').Value = '''2.318.13'
$ws.Range('E18').Value = '  +1.74%  '
$ws.Range('D19').Value = '''10.54'
$ws.Range('E19').Value = '  +2.84%  '
$ws.Range('D20').Value = '''4.29'
$ws.Range('E20').Value = '  +3.75%  '
$ws.Range('D21').Value = '''321.52'
$ws.Range('E21').Value = '  +6.55%  '
$ws.Range('D22').Value = '''6.67'
$ws.Range('E22').Value = '  +5.36%  '
$ws.Range('D23').Value = '''1.00'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').Value = '''60.71'
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('D25').Value = '''0.996'
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').Value = '''0.159'
$ws.Range('E26').Value = '  +6.51%  '
$ws.Range('D27').Value = '''7.66'
$ws.Range('E27').Value = '  +4.49%  '
$ws.Range('D28').Value = '''171.06'
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').Value = '''1.72'
$ws.Range('E29').Value = '  +7.55%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '''0.0₃0734'
$ws.Range('E30').Value = '  +6.25%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = '''1.20'
$ws.Range('E31').Value = '  +10.28%  '
$ws.Range('D32').Value = '''6.27'
$ws.Range('E32').Value = '  +5.26%  '
$ws.Range('D33').Value = '''18.34'
$ws.Range('E33').Value = '  +2.99%  '
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D35').Value = '''0.992'
$ws.Range('E35').Value = '  -0.43%  '
$ws.Range('D36').Value = '''1.27'
$ws.Range('E36').Value = '  +6.12%  '
$ws.Range('D37').Value = '''0.927'
$ws.Range('E37').Value = '  -1.62%  '
$ws.Range('D38').Value = '''4.01'
$ws.Range('E38').Value = '  +8.14%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '''1.52'
$ws.Range('E39').Value = '  +8.57%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = '''37.45'
$ws.Range('E40').Value = '  +4.11%  '
$ws.Range('D41').Value = '''0.383'
$ws.Range('E41').Value = '  +2.30%  '
$ws.Range('D42').Value = '''139.50'
$ws.Range('E42').Value = '  +11.48%  '
$ws.Range('D43').Value = '''3.58'
$ws.Range('E43').Value = '  +6.18%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '''5.17'
$ws.Range('E44').Value = '  +7.58%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').Value = '''275.56'
$ws.Range('E45').Value = '  +14.04%  '
$ws.Range('D46').Value = '''0.0510'
$ws.Range('E46').Value = '  +3.50%  '
$ws.Range('D47').Value = '''0.0930'
$ws.Range('E47').Value = '  +3.84%  '
$ws.Range('D48').Value = '''0.557'
$ws.Range('E48').Value = '  +1.81%  '
$ws.Range('E49').Value = '  +2.47%  '
$ws.Range('D50').Value = '''0.0216'
$ws.Range('E50').Value = '  +5.47%  '
$ws.Range('D51').Value = '''16.92'
$ws.Range('E51').Value = '  +4.95%  '
